# Fix the figure-label typos in the "illustrations" slide:
#   "A- Metabolic "  -> "A - " / "Metabolic "   (two runs, same formatting)
#   "B- Chemical "   -> "B - " / "Chemical "    (two runs, same formatting)
#   "B- "            -> "C "  / "- "            (two runs, same formatting)
# Each textbox is also nudged/resized slightly (spAutoFit reflow after the
# text change) per the target xfrm values below.

$p = $ppt.ActivePresentation

function Find-ShapeByName {
    param($slide, [string]$name)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

# PowerPoint's Shape.Left/.Top/.Width/.Height are in points (1 pt = 12700 EMU)
# while the OOXML stores integer EMUs. A tiny epsilon nudge keeps the
# point->EMU round trip from truncating down by one EMU.
function Emu-ToPoints {
    param([double]$emu)
    return ($emu / 12700.0) + 0.00004
}

# Find the slide that holds the figure (search every slide, don't assume index).
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $candidate = $p.Slides.Item($si)
    if ((Find-ShapeByName $candidate "TextBox 107") -ne $null) {
        $targetSlide = $candidate
        break
    }
}
$s = $targetSlide

# ---------------------------------------------------------------------
# Shape "TextBox 107": "A- Metabolic labelling" -> "A - Metabolic labelling"
# ---------------------------------------------------------------------
$shpA = Find-ShapeByName $s "TextBox 107"
$shpA.Left = Emu-ToPoints 325183
$shpA.Width = Emu-ToPoints 1621662

$trA = $shpA.TextFrame.TextRange
$trA.Characters(1, 3).Text = "A - "

# ---------------------------------------------------------------------
# Shape "TextBox 108": "B- Chemical labelling" -> "B - Chemical labelling"
# ---------------------------------------------------------------------
$shpB = Find-ShapeByName $s "TextBox 108"
$shpB.Left = Emu-ToPoints 2725753
$shpB.Width = Emu-ToPoints 1552733

$trB = $shpB.TextFrame.TextRange
$trB.Characters(1, 3).Text = "B - "

# ---------------------------------------------------------------------
# Shape "TextBox 207": "B- Label free" -> "C - Label free"
# Here the entire first run is "B- " (no trailing text sharing the run),
# so a plain text replacement would collapse back into a single run.
# Nudging the Bold flag (already True) on the first half forces PowerPoint
# to keep the edited span as its own run, matching the target two-run split.
# ---------------------------------------------------------------------
$shpC = Find-ShapeByName $s "TextBox 207"
$shpC.Left = Emu-ToPoints 5193257
$shpC.Width = Emu-ToPoints 1013291

$trC = $shpC.TextFrame.TextRange
$trC.Characters(1, 3).Text = "C - "
$trC.Characters(1, 2).Font.Bold = $true
